$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 68; existing rows 68..77 shift down to 69..78
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly price record
$ws.Cells.Item(68, 1).Value2 = 4
$ws.Cells.Item(68, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(68, 3).Value2 = "Los Lagos"
$ws.Cells.Item(68, 4).Value2 = 45034
$ws.Cells.Item(68, 5).Value2 = 10
$ws.Cells.Item(68, 6).Value2 = 100112030
$ws.Cells.Item(68, 7).Value2 = "Poroto granado"
$ws.Cells.Item(68, 8).Value2 = "Sin especificar"
$ws.Cells.Item(68, 9).Value2 = "Primera"
$ws.Cells.Item(68, 10).Value2 = 50
$ws.Cells.Item(68, 11).Value2 = 45000
$ws.Cells.Item(68, 12).Value2 = 45000
$ws.Cells.Item(68, 13).Value2 = 45000
$ws.Cells.Item(68, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(68, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(68, 16).Value2 = 1800
$ws.Cells.Item(68, 17).Value2 = 25
$ws.Cells.Item(68, 18).Value2 = "Hortaliza"
